$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Move Robot40 to location (4, 8) and remove the toolkit."
$ws.Range("B2").Value = $false
$ws.Range("F2").Value = $false

# Row 3
$ws.Range("A3").Value = "Move Robot40 to location (6, 2) and remove the liquid spill."
$ws.Range("B3").Value = $false
$ws.Range("C3").Value = $false
$ws.Range("D3").Value = $false
$ws.Range("E3").Value = $false
$ws.Range("F3").Value = $false
$ws.Range("G3").Value = $false

# Row 4
$ws.Range("A4").Value = "Move Robot9 to location (12, 3) and remove the large debris."
$ws.Range("B4").Value = $true

# Row 5
$ws.Range("A5").Value = "Move Robot35 to location (2, 11) and remove the dust."
$ws.Range("B5").Value = $false
$ws.Range("F5").Value = $false

# Row 6
$ws.Range("A6").Value = "Move Robot26 to location (12, 1) and remove the grass."
$ws.Range("B6").Value = $true

# Row 7
$ws.Range("A7").Value = "Move Robot41 to location (4, 11) and remove the small debris."

# Row 8
$ws.Range("A8").Value = "Move Robot2 to location (9, 1) and remove the vehicle."

# Row 9
$ws.Range("A9").Value = "Move Robot28 to location (11, 6) and remove the construction materials."

# Row 10
$ws.Range("A10").Value = "Move Robot39 to location (6, 1) and remove the tree branches."

# Row 11
$ws.Range("A11").Value = "Move Robot28 to location (2, 8) and remove the screws."
